$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-computed NATMI TPM edge table: updated specificity/weight values for the
# existing FAPs-sending rows, plus four new rows for the Resolving-Mac sending
# cluster (rows 6-9) that the refreshed pipeline run now emits.
$rows = @(
  @("FAPs", "Il21", "Il2rg", "ECs", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 35.32109533333334, 105.963286, 0.4123245624288747, 0.4123245624288747, 2.430974386316667, 21.87876947685, 0.2494176203586572, 0.2494176203586572),
  @("FAPs", "Il21", "Il2rg", "FAPs", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 1.228158333333333, 3.684475, 0.01433703690686912, 0.01433703690686912, 0.08452799729166666, 0.760751975625, 0.00867256029386408, 0.008672560293864078),
  @("FAPs", "Il21", "Il2rg", "MuSCs", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 2.583168333333334, 7.749505, 0.03015489023401347, 0.03015489023401347, 0.1777865605416667, 1.600079044875, 0.01824087539204396, 0.01824087539204395),
  @("FAPs", "Il21", "Il2rg", "Resolving-Mac", 1, 0.3333333333333333, 0.068825, 0.206475, 0.6049060451053807, 0.6049060451053806, 3, 1, 46.53090866666667, 139.592726, 0.5431835104302428, 0.5431835104302427, 3.202489788983333, 28.82240810085, 0.3285749890608155, 0.3285749890608154),
  @("Resolving-Mac", "Il21", "Il2rg", "ECs", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 35.32109533333334, 105.963286, 0.4123245624288747, 0.4123245624288747, 1.587789198519334, 14.290102786674, 0.1629069420702175, 0.1629069420702175),
  @("Resolving-Mac", "Il21", "Il2rg", "FAPs", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 1.228158333333333, 3.684475, 0.01433703690686912, 0.01433703690686912, 0.05520940155833333, 0.496884614025, 0.005664476613005042, 0.005664476613005041),
  @("Resolving-Mac", "Il21", "Il2rg", "MuSCs", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 2.583168333333334, 7.749505, 0.03015489023401347, 0.03015489023401347, 0.1161211660883333, 1.045090494795, 0.01191401484196952, 0.01191401484196952),
  @("Resolving-Mac", "Il21", "Il2rg", "Resolving-Mac", 1, 0.3333333333333333, 0.044953, 0.134859, 0.3950939548946194, 0.3950939548946193, 3, 1, 46.53090866666667, 139.592726, 0.5431835104302428, 0.5431835104302427, 2.091703937292667, 18.825335435634, 0.2146085213694274, 0.2146085213694273)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $i + 2
  $vals = $rows[$i]
  for ($j = 0; $j -lt $vals.Count; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $vals[$j]
  }
}
